# "update field beauty - MDLWL"
# The MDL record's id/claim-number field in row 2 (column A) is updated to a
# new value. (The workbook's shared-string pool also gains the other
# freshly-generated candidate values from the same data source, but those
# aren't referenced by any cell and the engine's own writer prunes unused
# shared strings on save, so only the live cell value matters here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CA-DNLQDJ5F"
